# Update the "想去人数" (Want-to-go count) column F values across sheets
# as produced by the gh-pages data generation run at commit 456a3b4.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 0
$ws.Cells.Item(3, 6).Value = 0
$ws.Cells.Item(4, 6).Value = 14560
$ws.Cells.Item(5, 6).Value = 0
$ws.Cells.Item(7, 6).Value = 0
$ws.Cells.Item(8, 6).Value = 53
$ws.Cells.Item(9, 6).Value = 0
$ws.Cells.Item(11, 6).Value = 55
$ws.Cells.Item(14, 6).Value = 0
$ws.Cells.Item(15, 6).Value = 41
$ws.Cells.Item(16, 6).Value = 30
$ws.Cells.Item(17, 6).Value = 136
$ws.Cells.Item(19, 6).Value = 0
$ws.Cells.Item(20, 6).Value = 0
$ws.Cells.Item(21, 6).Value = 76
$ws.Cells.Item(22, 6).Value = 0
$ws.Cells.Item(24, 6).Value = 7193
$ws.Cells.Item(25, 6).Value = 0
$ws.Cells.Item(26, 6).Value = 0
$ws.Cells.Item(27, 6).Value = 0
$ws.Cells.Item(28, 6).Value = 1167
$ws.Cells.Item(29, 6).Value = 0
$ws.Cells.Item(30, 6).Value = 0
$ws.Cells.Item(31, 6).Value = 0
$ws.Cells.Item(32, 6).Value = 42
$ws.Cells.Item(33, 6).Value = 0
$ws.Cells.Item(35, 6).Value = 220
$ws.Cells.Item(36, 6).Value = 5065
$ws.Cells.Item(37, 6).Value = 0

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(2, 6).Value = 0

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 0
$ws.Cells.Item(3, 6).Value = 0
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(5, 6).Value = 0
$ws.Cells.Item(7, 6).Value = 0
$ws.Cells.Item(8, 6).Value = 53
$ws.Cells.Item(9, 6).Value = 0
$ws.Cells.Item(11, 6).Value = 55
$ws.Cells.Item(14, 6).Value = 0
$ws.Cells.Item(15, 6).Value = 0
$ws.Cells.Item(16, 6).Value = 0
$ws.Cells.Item(17, 6).Value = 136
$ws.Cells.Item(19, 6).Value = 1320
$ws.Cells.Item(20, 6).Value = 0
$ws.Cells.Item(21, 6).Value = 76
$ws.Cells.Item(24, 6).Value = 0
$ws.Cells.Item(25, 6).Value = 7193
$ws.Cells.Item(27, 6).Value = 0
$ws.Cells.Item(28, 6).Value = 0
$ws.Cells.Item(29, 6).Value = 1167
$ws.Cells.Item(32, 6).Value = 0
$ws.Cells.Item(33, 6).Value = 0
$ws.Cells.Item(34, 6).Value = 42
$ws.Cells.Item(35, 6).Value = 0
$ws.Cells.Item(37, 6).Value = 0
$ws.Cells.Item(38, 6).Value = 0
$ws.Cells.Item(39, 6).Value = 0

